$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price cells are stored as literal text (e.g. "322.34"), not
# numbers, in the source workbook. Plain `.Value = "322.34"` assignment lets
# Excel's COM layer auto-coerce numeric-looking strings into real numbers
# (dropping trailing zeros, e.g. "0.0830" -> 0.083). Force text storage by
# stamping a text NumberFormat before the write, then restore the cell to
# its original, unformatted "Normal" style so no stray style index is left
# behind in the saved file.
function Set-TextValue {
    param($addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "43.740.66"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "2.245.26"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.16%  "
Set-TextValue "D5" "322.34"
$ws.Range("E5").Value = "  +2.35%  "
Set-TextValue "D6" "101.71"
$ws.Range("E6").Value = "  +0.72%  "
Set-TextValue "D7" "0.579"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -0.85%  "
Set-TextValue "D10" "37.41"
$ws.Range("E10").Value = "  +1.05%  "
Set-TextValue "D11" "0.0830"
$ws.Range("E11").Value = "  +0.32%  "
Set-TextValue "D12" "7.71"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "2.586.50"
$ws.Range("E14").Value = "  +0.96%  "
Set-TextValue "D15" "0.857"
$ws.Range("E15").Value = "  -0.30%  "
Set-TextValue "D16" "14.20"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "2.240.26"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").Value = "43.618.85"
$ws.Range("E18").Value = "  +1.49%  "
Set-TextValue "D19" "13.77"
$ws.Range("E19").Value = "  -4.59%  "
$ws.Range("D20").Value = "0.0₃0986"
$ws.Range("E20").Value = "  +2.34%  "
Set-TextValue "D21" "6.46"
$ws.Range("E21").Value = "  -0.60%  "
Set-TextValue "D22" "65.35"
$ws.Range("E22").Value = "  -0.43%  "
Set-TextValue "D23" "3.17"
$ws.Range("E23").Value = "  +1.37%  "
Set-TextValue "D24" "236.80"
$ws.Range("E24").Value = "  -0.37%  "
Set-TextValue "D25" "2.16"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("E26").Value = "  +0.15%  "
Set-TextValue "D27" "10.10"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("E28").Value = "  -1.49%  "
Set-TextValue "D29" "37.19"
$ws.Range("E29").Value = "  +8.15%  "
Set-TextValue "D30" "6.29"
$ws.Range("E30").Value = "  -2.06%  "
Set-TextValue "D31" "160.29"
$ws.Range("E31").Value = "  +4.09%  "
$ws.Range("E32").Value = "  -1.23%  "
Set-TextValue "D33" "0.0853"
$ws.Range("E33").Value = "  -2.39%  "
Set-TextValue "D34" "2.70"
$ws.Range("E34").Value = "  -2.01%  "
Set-TextValue "D35" "3.19"
$ws.Range("E35").Value = "  +4.10%  "
Set-TextValue "D36" "0.114"
$ws.Range("E36").Value = "  +8.52%  "
Set-TextValue "D37" "1.93"
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("E38").Value = "  -1.51%  "
Set-TextValue "D39" "3.80"
$ws.Range("E39").Value = "  +2.85%  "
Set-TextValue "D40" "4.30"
$ws.Range("E40").Value = "  -2.57%  "
Set-TextValue "D41" "15.35"
$ws.Range("E41").Value = "  +21.63%  "
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").Value = "1.811.67"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("E45").Value = "  -1.77%  "
Set-TextValue "D46" "82.91"
$ws.Range("E46").Value = "  -3.96%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D47" "5.25"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D48" "1.72"
$ws.Range("E48").Value = "  +5.97%  "
Set-TextValue "D49" "74.50"
$ws.Range("E49").Value = "  -3.59%  "
Set-TextValue "D50" "59.04"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("E51").Value = "  +0.52%  "